$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the split runs in the four "Завдання" list items back into single
#    runs (the wording itself is unchanged, only the run boundaries move).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ознайомитись із особливостями виникнення і поширення похибок",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ознайомитись із особливостями виникнення і поширення похибок", 2) | Out-Null

$d.Content.Find.Execute(
    "Оцінити похибку обмеження при обчисленні функції розкладом у ряд для заданого варіанта",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Оцінити похибку обмеження при обчисленні функції розкладом у ряд для заданого варіанта", 2) | Out-Null

$d.Content.Find.Execute(
    "Оцінити похибку заокруглення для заданого варіанта",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Оцінити похибку заокруглення для заданого варіанта", 2) | Out-Null

$d.Content.Find.Execute(
    "Дослідити поширення похибок для заданого варіанта",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Дослідити поширення похибок для заданого варіанта", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark off the "Варіант №14" paragraph - it will be
#    re-added later, at the end of the document, once the tail paragraphs
#    have been cleaned up.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3) Remove the "Результат" / "Висновки" sections that trail the results
#    table, collapsing everything back down to a single empty paragraph.
# ---------------------------------------------------------------------------
$t = $d.Tables(1)
$tEnd = $t.Range.End

# paragraph: <w:ind w:firstLine="0"/> (empty)
$rp = $d.Range($tEnd, $d.Content.End)
$para = $rp.Paragraphs(1)
$para.Range.Delete()

# paragraph: Heading2 "Результат"
$rp = $d.Range($tEnd, $d.Content.End)
$para = $rp.Paragraphs(1)
$para.Range.Delete()

# the next paragraph (first empty Consolas/autoSpaceDE paragraph) survives -
# remember where it ends so we can keep targeting what follows it.
$rp = $d.Range($tEnd, $d.Content.End)
$survivor = $rp.Paragraphs(1)
$survivorEnd = $survivor.Range.End

# paragraph: duplicate empty Consolas/autoSpaceDE paragraph
$rp = $d.Range($survivorEnd, $d.Content.End)
$para = $rp.Paragraphs(1)
$para.Range.Delete()

# paragraph: Heading2 "Висновки"
$rp = $d.Range($survivorEnd, $d.Content.End)
$para = $rp.Paragraphs(1)
$para.Range.Delete()

# Final paragraph holds the "Ознайомився із механізмами ..." conclusion text.
# Clear its text, then merge its (now empty) paragraph mark with the
# surviving paragraph immediately before it so only one paragraph remains.
$rp = $d.Range($survivorEnd, $d.Content.End)
$lastPara = $rp.Paragraphs(1)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
if ($lastTextRange.Start -lt $lastTextRange.End) {
    $lastTextRange.Delete()
}

$markRange = $d.Range($survivorEnd - 1, $survivorEnd)
$markRange.Delete()

# ---------------------------------------------------------------------------
# 4) Re-insert the "_GoBack" bookmark into the now-single trailing paragraph.
# ---------------------------------------------------------------------------
$finalRange = $d.Range($tEnd, $tEnd)
$d.Bookmarks.Add("_GoBack", $finalRange) | Out-Null

# ---------------------------------------------------------------------------
# 5) The footer's cached PAGE field result needs to reflect the new page
#    count.
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers(1)
    $ftr.Range.Find.Execute("2", $false, $false, $false, $false, $false, $true, 1, $false, "3", 2) | Out-Null
}
